$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 4681
$ws1.Cells.Item(3, 6).Value = 1857
$ws1.Cells.Item(4, 6).Value = 143
$ws1.Cells.Item(6, 6).Value = 3141
$ws1.Cells.Item(7, 6).Value = 582
$ws1.Cells.Item(8, 6).Value = 593
$ws1.Cells.Item(9, 6).Value = 277
$ws1.Cells.Item(10, 6).Value = 639
$ws1.Cells.Item(11, 6).Value = 548
$ws1.Cells.Item(12, 6).Value = 540
$ws1.Cells.Item(13, 6).Value = 395
$ws1.Cells.Item(15, 6).Value = 1789
$ws1.Cells.Item(16, 6).Value = 1355
$ws1.Cells.Item(18, 6).Value = 1629
$ws1.Cells.Item(19, 6).Value = 17
$ws1.Cells.Item(23, 6).Value = 46
$ws1.Cells.Item(26, 6).Value = 55
$ws1.Cells.Item(27, 6).Value = 109
$ws1.Cells.Item(32, 6).Value = 3920
$ws1.Cells.Item(33, 6).Value = 6
$ws1.Cells.Item(34, 6).Value = 774
$ws1.Cells.Item(35, 6).Value = 79
$ws1.Cells.Item(36, 6).Value = 1204
$ws1.Cells.Item(37, 6).Value = 59
$ws1.Cells.Item(38, 6).Value = 1868

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 6).Value = 26
$ws2.Cells.Item(3, 6).Value = 51

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 4681
$ws4.Cells.Item(3, 6).Value = 1857
$ws4.Cells.Item(4, 6).Value = 143
$ws4.Cells.Item(6, 6).Value = 3141
$ws4.Cells.Item(7, 6).Value = 582
$ws4.Cells.Item(8, 6).Value = 593
$ws4.Cells.Item(9, 6).Value = 277
$ws4.Cells.Item(10, 6).Value = 639
$ws4.Cells.Item(11, 6).Value = 548
$ws4.Cells.Item(12, 6).Value = 540
$ws4.Cells.Item(13, 6).Value = 26
$ws4.Cells.Item(14, 6).Value = 395
$ws4.Cells.Item(16, 6).Value = 1789
$ws4.Cells.Item(17, 6).Value = 1356
$ws4.Cells.Item(19, 6).Value = 1629
$ws4.Cells.Item(20, 6).Value = 17
$ws4.Cells.Item(24, 6).Value = 46
$ws4.Cells.Item(27, 6).Value = 55
$ws4.Cells.Item(28, 6).Value = 109
$ws4.Cells.Item(33, 6).Value = 3920
$ws4.Cells.Item(34, 6).Value = 51
$ws4.Cells.Item(35, 6).Value = 6
$ws4.Cells.Item(37, 6).Value = 774
$ws4.Cells.Item(38, 6).Value = 79
$ws4.Cells.Item(39, 6).Value = 1204
$ws4.Cells.Item(40, 6).Value = 59
$ws4.Cells.Item(41, 6).Value = 1868
